# Applies the "prace na vraceni zmen, oprava corrupted obrazku v prohlizeci" edit.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "ip_address_list" (sheet1): was 5 data rows, becomes 3.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ip_address_list")

# Drop the old rows 4 ("Domaci Wifi") and 5 ("514_Teleflex") entirely.
$ws1.Rows.Item(5).Delete()
$ws1.Rows.Item(4).Delete()

# Row 1
$ws1.Range("A1").Value = "514_Teleflexgg"
$ws1.Range("B1").Value = "192.168.14.240a"
$ws1.Range("C1").Value = "255.255.255.0"
$ws1.Range("D1").Value = "PC:192.168.14.240d"
$ws1.Range("E1").Value = 0

# Row 2
$ws1.Range("A2").Value = "Df gga"
$ws1.Range("B2").Value = "192.168.1.131g"
$ws1.Range("C2").Value = "255.255.255.0"
$ws1.Range("D2").ClearContents()
$ws1.Range("E2").Value = $true

# Row 3
$ws1.Range("A3").Value = "Domaci Wifiaffz"
$ws1.Range("B3").Value = "192.168.1.13¨ks"
$ws1.Range("C3").Value = "255.255.255.0"
$ws1.Range("D3").Value = "ddassssaa"
$ws1.Range("E3").Value = $false

# ---------------------------------------------------------------------------
# Sheet "ip_adress_fav_list" (sheet2): was empty, gets 2 new rows.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ip_adress_fav_list")

$ws2.Range("A1").Value = "Df gg"
$ws2.Range("B1").Value = "192.168.1.131g"
$ws2.Range("C1").Value = "255.255.255.0"
$ws2.Range("E1").Value = 1

$ws2.Range("A2").Value = "Df gga"
$ws2.Range("B2").Value = "192.168.1.131g"
$ws2.Range("C2").Value = "255.255.255.0"
$ws2.Range("E2").Value = $true

# ---------------------------------------------------------------------------
# Sheet "disk_list" (sheet3): was 6 data rows, becomes 4; hyperlink + style removed.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("disk_list")

# Drop the hyperlink on C1 (and the relationship that backs it).
$ws3.Hyperlinks.Delete()

# Drop the old rows 5 ("518_Valeo") and 6 ("Domaci Nas") entirely.
$ws3.Rows.Item(6).Delete()
$ws3.Rows.Item(5).Delete()

# Row 1 (loses the Hyperlink-style formatting that used to sit on C1)
$ws3.Range("A1").Value = "515_ZF"
$ws3.Range("B1").Value = "Z"
$ws3.Range("C1").Style = "Normal"
$ws3.Range("C1").Value = "\\10.9.250.100\08_Project_ZF_515\kamery"
$ws3.Range("D1").Value = "jhvadmin"
$ws3.Range("E1").Value = "jhvadm1n"
$ws3.Range("F1").ClearContents()

# Row 2
$ws3.Range("A2").Value = "518_Valeo II"
$ws3.Range("B2").Value = "V"
$ws3.Range("C2").Value = "\\192.168.1.10\10_vision"
$ws3.Range("D2").Value = "jhv_vision"
$ws3.Range("E2").Value = "Jhv*2708"
$ws3.Range("F2").Value = "Druha sít, ixon"

# Row 3
$ws3.Range("A3").Value = "Domaci Nas"
$ws3.Range("B3").Value = "S"
$ws3.Range("C3").Value = "\\192.168.1.20\Data"
$ws3.Range("D3").ClearContents()
$ws3.Range("E3").ClearContents()
$ws3.Range("F3").ClearContents()

# Row 4
$ws3.Range("A4").Value = "514_Teleflex"
$ws3.Range("B4").Value = "T"
$ws3.Range("C4").Value = "\\192.168.14.245\Data\Kamery"
$ws3.Range("D4").Value = "Vision"
$ws3.Range("E4").Value = "*Jhv2708"
$ws3.Range("F4").ClearContents()

# ---------------------------------------------------------------------------
# Sheet "Settings" (sheet4): a few default-value tweaks.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Settings")
$ws4.Range("B1").Value = 6
$ws4.Range("B3").Value = 0
$ws4.Range("B4").Value = 0

# ---------------------------------------------------------------------------
# Sheet "projects_bin2" (sheet5, hidden): was empty, gets 1 new row.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("projects_bin2")
$ws5.Range("A1").Value = "Df  fsifa"
$ws5.Range("B1").Value = "192.168.1.131g"
$ws5.Range("C1").Value = "255.255.255.0"
$ws5.Range("E1").Value = 0
